$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Step 1: advance the price-list date in A1 by one day (45308 -> 45309)
$ws.Range("A1").Value = 45309

# Step 2: update the unit prices for the four hinge sizes
$ws.Range("D22").Value = 12264
$ws.Range("D23").Value = 13894
$ws.Range("D24").Value = 18098
$ws.Range("D25").Value = 20048
